# Changes of 26th july 2022
# Update row 5's tracking number (P5) and actual rate (Q5) with new values.
# A leading apostrophe forces the numeric-looking text to be stored as a
# string (matching the original cell type) instead of being auto-converted
# to a number/currency value; resetting the style back to "Normal"
# afterwards avoids leaving a stray number-format/quote-prefix style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P5").Value = "'320018747466"
$ws.Range("Q5").Value = "'`$48.39"

$ws.Range("P5").Style = "Normal"
$ws.Range("Q5").Style = "Normal"
